# Brooklyn 2023 roster update: replace Nerlens Noel with Moses Brown (row 18)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Player name
$ws.Range("C18").Value = "Moses Brown"

# Jersey number (No.) cell is no longer populated for this player - remove it entirely
$ws.Range("B18").ClearContents()

# Height
$ws.Range("E18").Value = "7-2"

# Weight
$ws.Range("F18").Value = 245

# Birth date (kept as text, same style as the other rows in this column)
$ws.Range("G18").Value = "October 13, 1999"

# Years of experience - stored as text in this column (like the rest of the
# column). A leading apostrophe tells Excel to keep the numeric-looking
# string as text instead of coercing it into a number.
$ws.Range("I18").Value = "'3"

# College
$ws.Range("J18").Value = "UCLA"

# bbref url (displayed text of the hyperlink cell; underlying hyperlink relationship
# is left untouched, matching the source edit)
$ws.Range("K18").Value = "https://www.basketball-reference.com/players/b/brownmo01.html"
